$wb = $excel.ActiveWorkbook

# --- Step 1: insert the new "2022-Q1" sheet, positioned right before "总计" ---
$totalSheetBefore = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheetBefore)
$q1.Name = "2022-Q1"

# Re-fetch fresh references by name: Worksheets.Add() can rebind/invalidate
# handles obtained before the insert, so look everything up again by name.
$template = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Item("2022-Q1")
$totalSheet = $wb.Worksheets.Item("总计")

# Seed the new sheet with the same header row / layout / styles as the most
# recent quarter sheet (2021-Q4 already uses the "基金规模" header wording).
$template.Range("A1:H2").Copy($q1.Range("A1"))

# Overwrite the data row with the new quarter's fund data. These columns are
# stored as text in the source data (e.g. "0.60" keeps its trailing zero), so
# force text storage while typing them in, then drop the number-format
# override again so the cells end up unstyled just like the template row.
$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("B2").Value = "513080"
$q1.Range("C2").Value = "华安法国CAC40ETF（QDII）"
$q1.Range("D2").Value = "0.60"
$q1.Range("E2").Value = "96.69"
$q1.Range("F2").Value = "4.75"
$q1.Range("G2").Value = "0.0285"
$q1.Range("B2:G2").ClearFormats()
$q1.Range("H2").Value = 6

# --- Step 2: add a new top data row to "总计" sheet for 2022-Q1 ---
$totalSheet.Rows.Item(2).Insert()

# Inserting a row copies the format of the row above into the blank row; put
# the data cells back to the unstyled look the other data rows use, and give
# the index cell (column A) the same style as the rest of column A.
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.03

# Re-sequence the index column (A) for the rows that shifted down
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
